# Apply cryptos list price/volume update (GitHub Actions scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'62.694.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.65%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'2.908.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.22%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = "'  +0.06%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'568.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.05%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'144.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.00%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("E7").Value = "'  +0.08%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("E8").Value = "'  -1.22%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'2.905.21"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.08%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'6.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.16%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("E11").Value = "'  +1.06%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.433"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.48%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("E13").Value = "'  +1.26%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'32.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.31%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("E15").Value = "'  +0.46%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'3.390.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.26%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'62.661.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.60%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("E18").Value = "'  +0.64%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'2.910.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.15%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'430.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.31%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'13.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.81%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'0.660"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.08%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'6.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.40%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'78.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.73%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'11.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.29%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'10.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.38%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("E27").Value = "'  -0.08%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("E28").Value = "'  -1.63%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'0.0000111"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +6.09%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'7.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.26%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("E31").Value = "'  -2.08%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'2.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.30%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("E33").Value = "'  +0.00%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("B34").Value = "'EthereumClassic"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'25.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.01%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("B35").Value = "'Hedera"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'0.106"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.01%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'0.956"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.50%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = "'5.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.74%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'2.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.57%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'48.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.05%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("E40").Value = "'  -3.87%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("E41").Value = "'  -1.81%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'41.15"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +5.99%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'8.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.94%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'0.268"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.30%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'2.716.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.96%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'0.0339"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.20%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'133.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.22%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'358.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +4.72%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("E49").Value = "'  +0.00%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'0.000220"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +15.56%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("E51").Value = "'  -0.44%  "
$ws.Range("E51").Style = "Normal"
